# GACPAQ pages.xlsx - "Country Status" sheet update
# Updates several status cells (translation/app-review/main-study pipeline
# tracker) with new values, and refreshes the sheet view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Status")
$ws.Activate()

# xlPasteFormats = -4122 ; copies only cell formatting (fill/font/border),
# used to stamp the correct status-color style on each target cell by
# pulling it from a stable, untouched donor cell that already carries it.
# NOTE: this interpreter's named-parameter binding (-Name value) is
# unreliable, so the helper below takes positional args only.
function Set-StatusCell($Target, $StyleDonor, $Text) {
    $ws.Range($StyleDonor).Copy() | Out-Null
    $ws.Range($Target).PasteSpecial(-4122) | Out-Null
    $ws.Range($Target).Value = $Text
}

# Stable donor cells (style untouched by this edit):
#   C2  -> "completed" fill/font style (s=53)
#   G3  -> "on-going"/alt fill style   (s=56)
#   H14 -> "awaiting data" fill style  (s=54)

# Row 2
Set-StatusCell "J2" "H14" "in review"

# Row 3
Set-StatusCell "J3" "H14" "in review"

# Row 4
Set-StatusCell "I4" "C2" "completed"
Set-StatusCell "J4" "H14" "in review"

# Row 5
Set-StatusCell "F5" "G3" "in review"

# Row 6
Set-StatusCell "G6" "C2" "completed"
Set-StatusCell "H6" "G3" "being revised"
Set-StatusCell "I6" "G3" "being revised"
Set-StatusCell "J6" "H14" "in review"

# Row 7
Set-StatusCell "F7" "C2" "completed"
Set-StatusCell "G7" "C2" "completed"
Set-StatusCell "H7" "G3" "awaiting  data"
Set-StatusCell "I7" "G3" "awaiting  data"

# Row 8
Set-StatusCell "F8" "C2" "completed"
Set-StatusCell "G8" "C2" "completed"
Set-StatusCell "H8" "G3" "awaiting  data"
Set-StatusCell "I8" "G3" "awaiting  data"

# Row 13
Set-StatusCell "J13" "C2" "deployed"
Set-StatusCell "K13" "G3" "on-going (2)"

# Row 18
Set-StatusCell "I18" "C2" "completed"
Set-StatusCell "J18" "G3" "preparing…"

# Row 19
Set-StatusCell "K19" "G3" "on-going (1)"

# Row 20
Set-StatusCell "F20" "C2" "completed"
Set-StatusCell "G20" "C2" "completed"
Set-StatusCell "H20" "G3" "data received"
Set-StatusCell "I20" "G3" "data received"

# Row 21
Set-StatusCell "G21" "C2" "completed"
Set-StatusCell "H21" "H14" "awaiting  data"
Set-StatusCell "I21" "H14" "awaiting  data"

# Row 22
Set-StatusCell "G22" "C2" "completed"
Set-StatusCell "H22" "H14" "awaiting  data"
Set-StatusCell "I22" "H14" "awaiting  data"

$ws.Application.CutCopyMode = $false

# View-state refresh: zoom level and the active selection on the frozen pane.
$excel.ActiveWindow.Zoom = 125
$ws.Range("A18:XFD18").Select() | Out-Null
